# Apply updated cryptocurrency price/volume data to sheet1 (Coin / Link / Price / Volume(1h))
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.701.49"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "3.447.09"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.03"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.22"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "8.00"
$ws.Range("E9").Value = "  +4.58%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  +2.58%  "
$ws.Range("D12").Value = "4.037.13"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.37"
$ws.Range("E14").Value = "  -5.03%  "
$ws.Range("D15").Value = "3.440.81"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "62.726.82"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.36"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.69"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.10"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.08"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.564"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.30"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "3.581.17"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.67"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.02"
$ws.Range("E30").Value = "  -2.54%  "
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.35"
$ws.Range("E33").Value = "  -3.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.21"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("E35").Value = "  +4.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.37"
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "32.22"
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "169.68"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").Value = "3.479.86"
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0777"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.787"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.68"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("E45").Value = "  -2.37%  "
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").Value = "2.573.48"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.93"
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.72"
$ws.Range("E50").Value = "  -3.36%  "
$ws.Range("E51").Value = "  -0.05%  "
